# Fix validate room number is a number
# Add three more rooms to the Rooms sheet (106, 107, 303), preserving the
# "Room Number" column as text (matching the existing 101/102/103/104 cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rooms")

# --- Row 6: Room 106, price 0, Single, Free = Yes ---
$ws.Cells.Item(6, 1).Formula = '="106"'
$ws.Cells.Item(6, 1).Copy()
$ws.Cells.Item(6, 1).PasteSpecial(-4163)
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = "Single"
$ws.Cells.Item(6, 4).Value = "Yes"

# --- Row 7: Room 107, price 33, Single, Free = No ---
$ws.Cells.Item(7, 1).Formula = '="107"'
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(7, 1).PasteSpecial(-4163)
$ws.Cells.Item(7, 2).Value = 33
$ws.Cells.Item(7, 3).Value = "Single"
$ws.Cells.Item(7, 4).Value = "No"

# --- Row 8: Room 303, price 0, Single, Free = Yes ---
$ws.Cells.Item(8, 1).Formula = '="303"'
$ws.Cells.Item(8, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4163)
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = "Single"
$ws.Cells.Item(8, 4).Value = "Yes"

$excel.CutCopyMode = 0
